$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report-year text in the "description" column (B) for all rows ---
# Each note now cites the specific report year (e.g. "113年性別統計年報" instead of
# the generic "性別統計年報"); row 3 additionally gets a distinct "112年" note.
$ws.Range("B2").Value = "113年工務局暨所屬職員總計524人，其中男生332人(占63.36%)、女生192人(占36.64%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B3").Value = "112年工務局暨所屬職員總計519人，其中男生331人(占63.78%)、女生188人(占36.22%)。`n（資料來源）`n高雄市政府工務局112年性別統計年報。"
$ws.Range("B4").Value = "113年工務局職員總計191人，其中男生118人(占61.78%)、女生70人(占38.22%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B5").Value = "113年工務局職員簡任6人，其中男生5人(占83.33%)、女生1人(占16.67%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B6").Value = "113年工務局職員薦任143人，其中男生94人(占65.73%)、女生49人(占34.27%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B7").Value = "113年工務局職員委任39人，其中男生19人(占48.72%)、女生20人(占51.28%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B8").Value = "113年新建工程處職員總計131人，其中男生80人(占61.07%)、女生51人(占38.93%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B9").Value = "113年新建工程處職員簡任3人，均為男生。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B10").Value = "113年新建工程處職員薦任94人，其中男生59人(占62.77%)、女生35人(占37.23%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B11").Value = "113年新建工程處職員委任34人，其中男生18人(占52.94%)、女生16人(占47.06%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B12").Value = "113年道路養護工程處職員總計82人，其中男生53人(占64.63%)、女生29人(占35.37%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B13").Value = "113年道路養護工程處職員簡任3人，均為男生。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B14").Value = "113年道路養護工程處職員薦任49人，其中男生34人(占69.39%)、女生15人(占30.61%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B15").Value = "113年道路養護工程處職員委任30人，其中男生16人(占53.33%)、女生14人(占46.67%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B16").Value = "113年公園處職員總計94人，其中男生59人(占62.77%)、女生35人(占37.23%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B17").Value = "113年公園處職員簡任3人，均為男生。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B18").Value = "113年公園處職員薦任61人，其中男生41人(占67.21%)、女生20人(占32.79%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B19").Value = "113年公園處職員委任30人，其中男生15人(占50.00%)、女生15人(占50.00%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B20").Value = "113年違章建築處理大隊職員總計30人，其中男生23人(占76.67%)、女生7人(占23.33%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B21").Value = "113年違章建築處理大隊職員薦任15人，其中男生13人(占86.67%)、女生2人(占13.33%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B22").Value = "113年違章建築處理大隊職員委任15人，其中男生10人(占66.67%)、女生5人(占33.33%)。`n（資料來源）`n高雄市政府工務局113年性別統計年報。"
$ws.Range("B23").Value = "113年工務局主管預算數合計8,194,228千元，執行數5,786,539千元，執行率89.70%。`n（資料來源）`n高雄市政府工務局113年工務統計年報。"
$ws.Range("B24").Value = "113年工務局主管經常門預算數2,023,521千元，執行數1,952,014千元，執行率96.47%。`n（資料來源）`n高雄市政府工務局113年工務統計年報。"
$ws.Range("B25").Value = "113年工務局主管資本門預算數6,170,707千元，執行數5,831,206千元，執行率94.50%。`n（資料來源）`n高雄市政府工務局113年工務統計年報。"
$ws.Range("B26").Value = "113年工務局主管預算數合計8,194,228千元，執行數5,786,539千元，執行率89.70%。`n（資料來源）`n高雄市政府工務局113年工務統計年報。"

# --- Update the sheet view: scroll the frozen-header view down and move the active cell ---
$ws.Range("A24").Select()
